# Add a new "Heatmap setting" bullet item to the Settings File list,
# right after the "Storage type (1 = xls, 2 = xlsx, 3 = csv)" item and
# before the "Procedures" heading.

$d = $word.ActiveDocument

# Locate the "Storage type ..." list paragraph by scanning the
# paragraph collection for its text (robust against Find() collapsing
# the range to just the matched text instead of the whole paragraph).
$storageTypePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Storage type*") {
        $storageTypePara = $candidate
        break
    }
}

# Insert a brand-new paragraph right after it; Word clones the
# paragraph formatting (ListParagraph style / numId=2 list numbering)
# from the paragraph it follows.
$storageTypePara.Range.InsertParagraphAfter()
$p1 = $storageTypePara.Next()

# Populate the new paragraph with the first chunk of text.
$p1.Range.Text = "Heatmap "

# Build the remaining two chunks as their own paragraphs (so each gets
# its own run), then splice them back onto the first paragraph by
# deleting the paragraph marks between them. Removing a paragraph mark
# merges the two paragraphs' runs together without coalescing the text
# into a single run.
$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "setting"

$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$p3.Range.Text = " (A for automatic else provide a number to compare sets of heatmaps)"

$mark1 = $d.Range($p1.Range.End - 1, $p1.Range.End)
$mark1.Delete()

$p1again = $storageTypePara.Next()
$mark2 = $d.Range($p1again.Range.End - 1, $p1again.Range.End)
$mark2.Delete()
